# "Data wrangling" rewrite of Sheet1: turn the wide "Year" table
# (A1 "Year" / A2,A4 "5.442","5.474" header column + 12 side-by-side
# metric columns) into a long/melted layout that mirrors a pandas
# `df.to_excel()` dump:
#   - Row 1: 5-column header "Unnamed: 0", 2019, "Unnamed: 1", 2018,
#            "Unnamed: 2"
#   - Column B: the 2019 figures, Column D: the 2018 figures
#   - Row 5 intentionally left blank (mirrors a blank line in the
#     source dataframe), spacer column C/E left blank too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stash two "clean" format swatches from cells outside the region
# we are about to rewrite, so we can restore proper styling afterwards
# without minting a pile of brand-new style-pool entries:
#   Z1 <- A1's existing bold/centered/bordered "header" cell style
#   Z2 <- an untouched cell's plain default ("Normal") style
$ws.Range("A1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("Z2").Copy()
$excel.CutCopyMode = $false

# Wipe the old table.
$ws.Range("A1:M5").Clear()

# ---- Header row -------------------------------------------------------
$ws.Range("A1").Value = "Unnamed: 0"
$b1 = $ws.Range("B1")
$b1.NumberFormat = "@"   # "2019" looks numeric; force it to text first
$b1.Value = "2019"
$ws.Range("C1").Value = "Unnamed: 1"
$d1 = $ws.Range("D1")
$d1.NumberFormat = "@"   # "2018" looks numeric; force it to text first
$d1.Value = "2018"
$ws.Range("E1").Value = "Unnamed: 2"

# ---- Data rows (column B = 2019 values, column D = 2018 values) -------
# Row 5 is intentionally skipped. All values are written as text
# (e.g. "3,145", "5.442") rather than numbers, matching the source.
$values2019 = @("5.442", "3,145", "2,297", $null, "485", "249", "236", "2,685", "1,493", "1,192", "2,272", "1,403", "869")
$values2018 = @("5.474", "3,157", "2,317", $null, "471", "229", "242", "2,825", "1,594", "1,231", "2,178", "1,334", "844")

for ($i = 0; $i -lt $values2019.Length; $i++) {
    $row = $i + 2
    if ($null -ne $values2019[$i]) {
        $cell = $ws.Cells.Item($row, 2)
        $cell.NumberFormat = "@"
        $cell.Value = $values2019[$i]
    }
    if ($null -ne $values2018[$i]) {
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $values2018[$i]
    }
}

# --- Re-apply clean formatting from the stashed swatches: header row
# gets the bold/bordered style back, data cells go back to the plain
# default style (undoing the transient "@" text format above).
$ws.Range("Z1").Copy()
$ws.Range("A1:E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("Z2").Copy()
$ws.Range("B2:B4").PasteSpecial(-4122)
$ws.Range("D2:D4").PasteSpecial(-4122)
$ws.Range("B6:B14").PasteSpecial(-4122)
$ws.Range("D6:D14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Drop the temporary swatch cells.
$ws.Range("Z1:Z2").Clear()
